# "add to cart tc" — renumber the later TC sheets by +1 (a new TC sheet was
# inserted earlier in the suite), and move the active/selected tab from the
# old "EditDeleteAddress" sheet to the (renumbered) "AddToCart" sheet, whose
# own selection moves from E14 to E5.

$wb = $excel.ActiveWorkbook

# Rename sheets 6-11 (1-indexed) to bump their TC numbers by one.
$wb.Worksheets.Item(6).Name  = "TC06_EditProfile"
$wb.Worksheets.Item(7).Name  = "TC07_ChangePassword"
$wb.Worksheets.Item(8).Name  = "TC08_AddAddress"
$wb.Worksheets.Item(9).Name  = "TC09_EditDeleteAddress"
$wb.Worksheets.Item(10).Name = "TC10_AddToWishlist"
$wb.Worksheets.Item(11).Name = "TC11_AddToCart"

# The previously-active sheet (now "TC09_EditDeleteAddress") keeps its
# selection but loses tab focus; its active cell moves from B3 to B2.
$wsPrev = $wb.Worksheets.Item(9)
$wsPrev.Range("B2").Select()

# The new active sheet is the renumbered "TC11_AddToCart"; its active cell
# moves from E14 to E5.
$wsNew = $wb.Worksheets.Item(11)
$wsNew.Activate()
$wsNew.Range("E5").Select()
